$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Done?" column for the two existing PAGE RENDEZ-VOUS rows flips from YES to NO.
$ws.Range("C20").Value = "NO"
$ws.Range("C21").Value = "NO"

# The two bottom rows (25/26) get their previously-empty "Done?" cell filled in.
$ws.Range("C25").Value = "NO"
$ws.Range("C26").Value = "NO"

# Leave the view scrolled/selected where the author ended up when saving.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G31").Select()
